$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'small n'
$ws.Range("A2").Value = 'small rainbow'
$ws.Range("A3").Value = 'small soccer'
$ws.Range("A4").Value = 'small waist pack for women'
$ws.Range("A5").Value = 'smart basketball ball'
$ws.Range("A6").Value = 'smart basketball training'
$ws.Range("A7").Value = 'smart soccee ball'
$ws.Range("A8").Value = 'smart soccer ball size 4'
$ws.Range("A9").Value = 'smart soccor ball'
$ws.Range("A10").Value = 'smart white bread'
$ws.Range("A11").Value = 'snag tights'
$ws.Range("A12").Value = 'snowboarding compression pants'
$ws.Range("A13").Value = 'snowboarding gear for women'
$ws.Range("A14").Value = 'snowboarding gear women'
$ws.Range("A15").Value = 'snowboarding knee brace'
$ws.Range("A16").Value = 'snowboarding leggings women'
$ws.Range("A17").Value = 'snowboarding pants woman'
$ws.Range("A18").Value = 'so capri leggings'
$ws.Range("A19").Value = 'so yoga leggings'
$ws.Range("A20").Value = 'soft brushed leggings'
$ws.Range("A21").Value = 'soft comfy leggings'
$ws.Range("A22").Value = 'soft women leggings'
$ws.Range("A23").Value = 'sore hip'
$ws.Range("A24").Value = 'sore hip joint'
$ws.Range("A25").Value = 'sore hips'
$ws.Range("A26").Value = 'sore hips and legs'
$ws.Range("A27").Value = 'sore joints'
$ws.Range("A28").Value = 'sore knee support'
$ws.Range("A29").Value = 'sore legs'
$ws.Range("A30").Value = 'sore legs and hips'
$ws.Range("A31").Value = 'sore muscles and joints'
$ws.Range("A32").Value = 'sound skin roll'
$ws.Range("A33").Value = 'spandex capri leggings'
$ws.Range("A34").Value = 'spandex capris'
$ws.Range("A35").Value = 'spandex nike'
$ws.Range("A36").Value = 'spandex nike pro'
$ws.Range("A37").Value = 'spandex running pants women'
$ws.Range("A38").Value = 'spandex tank tops for women compression'
$ws.Range("A39").Value = 'spandex tights'
$ws.Range("A40").Value = 'spandex workout capri'
$ws.Range("A41").Value = 'spandex yoga pants'
$ws.Range("A42").Value = 'spans tights'
$ws.Range("A43").Value = 'spanx capri length'
$ws.Range("A44").Value = 'spanx clothing for women'
$ws.Range("A45").Value = 'spanx compression'
$ws.Range("A46").Value = 'spanx legging'
$ws.Range("A47").Value = 'spanx leggings'
$ws.Range("A48").Value = 'spanx medium'
$ws.Range("A49").Value = 'spanx men'
$ws.Range("A50").Value = 'spanx nike'
$ws.Range("A51").Value = 'spanx pants'
$ws.Range("A52").Value = 'spanx pants for women'
$ws.Range("A53").Value = 'spanx plus'
$ws.Range("A54").Value = 'spanx wear'
$ws.Range("A55").Value = 'spanxs tights'
$ws.Range("A56").Value = 'spartan apparel for women'
$ws.Range("A57").Value = 'spartan compression pants'
$ws.Range("A58").Value = 'spartan cross'
$ws.Range("A59").Value = 'spartan kinesiology tape'
$ws.Range("A60").Value = 'spartan race apparel for women'
$ws.Range("A61").Value = 'spartan race clothes for women'
$ws.Range("A62").Value = 'spartan race clothing'
$ws.Range("A63").Value = 'spartan race gear'
$ws.Range("A64").Value = 'spartan race gear men'
$ws.Range("A65").Value = 'spartan race pants'
$ws.Range("A66").Value = 'spartan race women'
$ws.Range("A67").Value = 'spartan race womens clothing'
$ws.Range("A68").Value = 'spartan training equipment'
$ws.Range("A69").Value = 'spartan workout gear'
$ws.Range("A70").Value = 'spectrum basketball'
$ws.Range("A71").Value = 'spectrum bread basket'
$ws.Range("A72").Value = 'spectrum net'
$ws.Range("A73").Value = 'spectrum optimum'
$ws.Range("A74").Value = 'speed gear'
$ws.Range("A75").Value = 'speed shorts womens'
$ws.Range("A76").Value = 'speed track shorts'
$ws.Range("A77").Value = 'spf 50 pants for women'
$ws.Range("A78").Value = 'spf leggings for women'
$ws.Range("A79").Value = 'spf pants women'
$ws.Range("A80").Value = 'spine of god'
$ws.Range("A81").Value = 'spine support for women'
$ws.Range("A82").Value = 'sport 4'
$ws.Range("A83").Value = 'sport bike pants'
$ws.Range("A84").Value = 'sport capri'
$ws.Range("A85").Value = 'sport capri pants for women'
$ws.Range("A86").Value = 'sport capri pants women'
$ws.Range("A87").Value = 'sport capris for women'
$ws.Range("A88").Value = 'sport compression'
$ws.Range("A89").Value = 'sport compression tights'
$ws.Range("A90").Value = 'sport hoodies women'
$ws.Range("A91").Value = 'sport joint'
$ws.Range("A92").Value = 'sport knee support'
$ws.Range("A93").Value = 'sport leggings'
$ws.Range("A94").Value = 'sport leggings plus size'
$ws.Range("A95").Value = 'sport leggings women'
$ws.Range("A96").Value = 'sport light'
$ws.Range("A97").Value = 'sport lycra fabric'
$ws.Range("A98").Value = 'sport pant'
$ws.Range("A99").Value = 'sport pants for women'
$ws.Range("A100").Value = 'sport pants for womens'
